$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1 (13:25 -> 13:50)
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 13:50"

# The province table (sorted descending by "Casos totales") had La Rioja's case
# count rise from 1436 to 1629, overtaking Zaragoza (1449). Swap the two rows so
# La Rioja now sits above Zaragoza, each carrying its updated figures.
$ws.Range("A13").Value = "La Rioja"
$ws.Range("B13").Value = 1629
$ws.Range("C13").Value = 397
$ws.Range("D13").Value = 1164
$ws.Range("E13").Value = 68

$ws.Range("A14").Value = "Zaragoza"
$ws.Range("B14").Value = 1449
$ws.Range("C14").Value = 101
$ws.Range("D14").Value = 1269
$ws.Range("E14").Value = 79
